# Refactor CU position bullets for clarity
# The three "Coordinated / Participated in the engineering / Acted as primary
# maintainer" bullets get rotated and re-worded, and the following three
# bullets (Puppet config mgmt, 3rd tier support, CUIT integration) get small
# clarifying edits. All edits are applied as whole/partial-sentence
# find & replace operations so run-level formatting (Tahoma, 10pt) is
# preserved automatically by Word.

$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false,
                                       $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $oldText"
    }
}

# Old bullet "Coordinated ..." becomes the new first bullet, built from the
# old "Acted as primary maintainer ..." bullet with wording tweaks.
Replace-Text "Coordinated within and across teams to plan and execute ongoing projects. Projects included Gitlab CI/CD for CAS deployments, VMware vSphere upgrades, SVN to Git Migration, and MediaWiki to Confluence migration." `
             "Acted as primary maintainer of various servers and platforms, such as CAS Single Sign On System (SSO), Azure DR environment, Drupal CMS, VMware vSphere, DHCP, Zenoss Monitoring, and others."

# Old bullet "Participated in the engineering ..." becomes the new second
# bullet, built from the old "Coordinated ..." bullet with wording tweaks.
Replace-Text "Participated in the engineering and evaluation of new installations and upgrades of hardware and software. Projects included staging server for CAS SSO in Azure, deployment of AWS CloudFront for static web assets, and vSphere upgrades." `
             "Coordinated within and across teams to plan and execute ongoing projects. Projects included Gitlab CI/CD for CAS deployments, additional VMware vSphere hosts, SVN to Git Migration, and MediaWiki to Confluence migration."

# Old bullet "Acted as primary maintainer ..." becomes the new third bullet,
# built from the old "Participated in the engineering ..." bullet with
# wording tweaks.
Replace-Text "Acted as primary maintainer of various web services/servers, such as CAS Single Sign On System (SSO), Azure DR environment, Drupal Systems, VMware vSphere, DHCP, Zenoss Monitoring, and others." `
             "Participated in the engineering and evaluation of new installations and upgrades of hardware and software. Projects included staging server for CAS SSO in Azure, deployment of AWS CloudFront for static web assets, and vSphere and vCenter upgrades."

# Puppet / GitLab code-approval bullet: reword the second sentence.
Replace-Text "Recently setup a code approval process in GitLab after the SVN to Git Migration of the Business School's Puppet code." `
             "Setup of a code approval process in GitLab for all changes."

# 3rd tier support bullet: clarify what kind of problems.
Replace-Text " tier support to the Client Support Group to resolve complex problems." `
             " tier support to the Client Support Group to resolve complex end user issues."

# CUIT integration bullet: reword network refresh / IDM systems clause.
Replace-Text "Access and Data Center network refresh to CUIT, and integration of IDM functions central Columbia University and Business School AD and IDM systems AD and IDM systems." `
             "Access and Data Center network re-architecture, and integration of IDM functions between central Columbia University and Business School AD, HRMS, and SSO systems."
